$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Introduction to Python")

# Update E4 and E6 remarks from "Try Again" to "Good"
$ws.Range("E4").Value = "Good"
$ws.Range("E6").Value = "Good"

# Add new column F with header and set its width
$ws.Range("F1").Value = "#Don't touch Medium Questions yet"
$ws.Range("F1").Font.Bold = $true
$ws.Columns.Item(6).ColumnWidth = 30.83

# Update the active selection
$ws.Range("A26").Select()

# Delete the "Daily Problems" worksheet
$ws2 = $wb.Worksheets.Item("Daily Problems")
$ws2.Delete()
